# Fill in the previously-missing "NA" measurements for the t1 samples
# (rows 2, 9, 16, 23) with the real numbers copied over from the
# measurement/R scripts, formatted with the small gray-blue Consolas font
# used in the source console output. AVERAGE/STDEV.S formulas in J/K
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (LIVER, sample t1) ---
$ws.Range("E2").Value = -364.89473700000002
$ws.Range("F2").Value = -479.78947399999998
$ws.Range("H2").Value = -604.89473699999996
$ws.Range("I2").Value = -305.90476200000001

# --- Row 9 (LIVER Gd ratio, sample t1) ---
$ws.Range("E9").Value = 0.37052600000000002
$ws.Range("F9").Value = 0.86421099999999995
$ws.Range("H9").Value = 0.148421
$ws.Range("I9").Value = 1.9

# --- Row 16 (SPLEEN, sample t1) ---
$ws.Range("E16").Value = -585.684211
$ws.Range("F16").Value = -422.57894700000003
$ws.Range("H16").Value = -652.95000000000005
$ws.Range("I16").Value = -122.38888900000001

# --- Row 23 (SPLEEN Gd ratio, sample t1) ---
$ws.Range("E23").Value = 0.24315800000000001
$ws.Range("F23").Value = 1.6926319999999999
$ws.Range("H23").Value = 0.223
$ws.Range("I23").Value = 2.605556

# Build the font (Consolas 8pt, light blue-gray FFD6DEEB, family "Modern")
# on the first touched cell, then paint that formatting onto the other
# newly-filled cells instead of redoing the same four Font assignments
# sixteen times.
$src = $ws.Range("E2")
$src.Font.Name = "Consolas"
$src.Font.Size = 8
$src.Font.Color = 15458006
$src.Font.Family = 3

$src.Copy()
$formatTargets = @("F2","H2","I2","E9","F9","H9","I9","E16","F16","H16","I16","E23","F23","H23","I23")
foreach ($addr in $formatTargets) {
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# Leave the selection where the author last left it before saving.
$ws.Range("I23").Select()
